# Sierra Leone master data: translate reg_center_type reference data from
# French (Madagascar) to English and refresh the sheet layout/formatting to
# match the re-exported workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column-level formatting first (before any row past 2 is touched, so
# whole-column writes don't materialize phantom cells in later rows). ------
$ws.Columns.Item(5).NumberFormat = "@"
$ws.Columns.Item(1).ColumnWidth = 8.43
$ws.Columns.Item(3).ColumnWidth = 21.07
$ws.Columns.Item(4).ColumnWidth = 33.26
$ws.Columns.Item(5).ColumnWidth = 7.62

# ---- Content edits -------------------------------------------------------
# Row 2 data: language + labels translated fra -> eng.
$ws.Range("A2").Value = "eng"
$ws.Range("C2").Value = "Regular"
$ws.Range("D2").Value = "Regular Registration Center"

# is_active was a boolean TRUE; it is now stored as literal text "TRUE".
$ws.Range("E2").Value = "'TRUE"

# ---- Header / data formatting refresh ------------------------------------
# Header row: keep bold/centered/top alignment but drop the wrap text.
$ws.Range("A1:E1").WrapText = $false
$ws.Range("A1:E1").HorizontalAlignment = -4108
$ws.Range("A1:E1").VerticalAlignment = -4160

# Data & header rows no longer carry the heavy medium borders.
$ws.Range("A1:E2").Borders.LineStyle = -4142

# Rows go back to their default auto height (no more thick-bottom/tall rows).
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()

# A later, still-empty row (data entry continued below row 2).
$ws.Range("C4:D4").WrapText = $true
$ws.Range("C4:D4").HorizontalAlignment = -4131

# ---- View / selection ------------------------------------------------
$ws.Range("E8").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100

# ---- Page setup refreshed to the "Normal"-ish margins + A4 portrait. -----
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36.85
$ws.PageSetup.FooterMargin = 36.85
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "reg_center_type translated to English master data"
